$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("huckleberry")

# Unhide columns A and B and set their widths (previously hidden with width 0)
$ws1.Columns.Item(1).Hidden = $false
$ws1.Columns.Item(1).ColumnWidth = 14.498697916666666
$ws1.Columns.Item(2).Hidden = $false
$ws1.Columns.Item(2).ColumnWidth = 21.053385416666668

# Make "huckleberry" the active/selected sheet (was "Formatted")
$ws1.Activate()
